# Add monthly load visualization:
#  - Sheet2: add "# of days (monthly)" / "# of days (quarterly)" columns (C, D)
#            with ROUND() formulas and a totals row.
#  - Sheet3 (new): day-by-day (1-30) load_profile lookup table, sorted by day.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet2 : add the day-count columns
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("C1").Value = "# of days (monthly)"
$ws2.Range("D1").Value = "# of days (quarterly)"

$probabilities = @(0.22, 0.27, 0.19, 0.17, 0.1, 0.04, 0.01)
for ($i = 0; $i -lt $probabilities.Count; $i++) {
    $row = 2 + $i
    $ws2.Range("B$row").Value = $probabilities[$i]
}

$ws2.Range("C2").Formula = "=ROUND(B2*30,0)"
$ws2.Range("D2").Formula = "=ROUND(B2*91.25,0)"
$ws2.Range("C3:C8").Formula = "=ROUND(B3*30,0)"
$ws2.Range("D3:D8").Formula = "=ROUND(B3*91.25,0)"

$ws2.Range("B9").Formula = "=SUM(B2:B8)"
$ws2.Range("C9").Formula = "=SUM(C2:C8)"
$ws2.Range("D9").Formula = "=SUM(D2:D8)"

$ws2.Columns.Item(3).ColumnWidth = 16.166666666666668
$ws2.Columns.Item(4).ColumnWidth = 17.166666666666668

$excel.Calculate()

$ws2.Range("B8").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet3 : new day -> load_profile lookup sheet
# (copy Sheet1 - which carries no extra column-width overrides - rather than
#  Worksheets.Add(), so the new sheet's default row height / column
#  formatting matches the rest of this workbook instead of the engine's
#  built-in blank-sheet defaults)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Copy([System.Reflection.Missing]::Value, $lastSheet)
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "Sheet3"
$ws3.Cells.Clear()

$ws3.Range("A1").Value = "day"
$ws3.Range("B1").Value = "load_profile"

$loadProfiles = @(
    "med_on_peak", "med_off_peak", "med_on_peak", "low_on_peak", "low_on_peak",
    "med_on_peak", "med_off_peak", "low_on_peak", "med_on_peak", "med_on_peak",
    "high_off_peak", "low_on_peak", "high_on_peak", "low_off_peak", "med_off_peak",
    "low_off_peak", "low_off_peak", "low_off_peak", "low_off_peak", "low_on_peak",
    "low_off_peak", "high_on_peak", "med_off_peak", "low_off_peak", "low_on_peak",
    "low_off_peak", "low_on_peak", "high_on_peak", "med_on_peak", "med_off_peak"
)

for ($i = 0; $i -lt $loadProfiles.Count; $i++) {
    $row = 2 + $i
    $ws3.Cells.Item($row, 1).Value = $i + 1
    $ws3.Cells.Item($row, 2).Value = $loadProfiles[$i]
}

$ws3.Columns.Item(2).ColumnWidth = 11.830729166666666

# Record the sort-by-day state (data is already in day order, but this
# mirrors the author applying Data > Sort on column A).
$sort = $ws3.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws3.Range("A2:A31")) | Out-Null
$sort.SetRange($ws3.Range("A2:B31"))
$sort.Header = 2
$sort.Apply()

$ws3.Activate() | Out-Null
$ws3.Range("D4").Select() | Out-Null
